$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.732.25"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.062.35"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'245.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("E6").Value = "  +1.01%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'55.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.65%  "

$ws.Range("D9").Value = "'60.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("E10").Value = "  -2.14%  "

$ws.Range("D11").Value = "'0.0753"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.73%  "

$ws.Range("E12").Value = "  -2.98%  "

$ws.Range("D13").Value = "'0.936"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.08%  "

$ws.Range("D14").Value = "'14.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.72%  "

$ws.Range("D15").Value = "2.364.15"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").Value = "'5.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.62%  "

$ws.Range("D17").Value = "2.069.27"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").Value = "36.676.16"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "'17.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.75%  "

$ws.Range("D20").Value = "'72.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "

$ws.Range("E21").Value = "  -2.25%  "

$ws.Range("D22").Value = "'238.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  -3.43%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  -2.44%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.72%  "

$ws.Range("D27").Value = "'9.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.38%  "

$ws.Range("D28").Value = "'166.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("E31").Value = "  +8.72%  "

$ws.Range("E32").Value = "  -6.76%  "

$ws.Range("E33").Value = "  -3.82%  "

$ws.Range("E34").Value = "  -2.58%  "

$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'0.0849"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.19%  "

$ws.Range("E38").Value = "  -2.90%  "

$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "'5.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.86%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.95%  "

$ws.Range("D41").Value = "'2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.75%  "

$ws.Range("D42").Value = "'0.0216"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.00%  "

$ws.Range("E43").Value = "  -4.11%  "

$ws.Range("D44").Value = "'95.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.74%  "

$ws.Range("D45").Value = "'0.0909"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.31%  "

$ws.Range("D46").Value = "1.413.58"
$ws.Range("E46").Value = "  +8.65%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.21%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'16.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.66%  "

$ws.Range("D49").Value = "'2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("D50").Value = "'2.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.82%  "

$ws.Range("D51").Value = "2.250.60"
$ws.Range("E51").Value = "  +0.39%  "
